# Screen-capture placeholders for the <diagnostics-general> section and the
# command-line "arg" spell-check split.
#
# Strategy: locate each target paragraph by its visible text with
# Find.Execute, then overwrite the *whole paragraph* (the paragraph's
# Range, which includes the trailing paragraph mark) via Range.InsertXML
# with the exact OOXML we want - this lets us control both the run-level
# <w:rPr> and the paragraph-mark <w:pPr><w:rPr> precisely, and also lets us
# insert the <w:proofErr/> spell-check markers that plain text/Font edits
# cannot produce.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# NOTE: `$range.Paragraphs(1)` on a narrow (Find-produced) sub-range is
# unreliable in this host - it can resolve to an unrelated paragraph
# elsewhere in the document instead of the one actually containing the
# range. Walk `$d.Paragraphs` (the whole-document collection, which does
# index correctly) and pick the paragraph whose extent contains our found
# range instead.
function Get-ContainingParagraph($rng) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            return $p
        }
    }
    throw "No containing paragraph found for range $($rng.Start)-$($rng.End)"
}

function Set-ParagraphXml($findText, $paragraphInnerXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $findText"
    }
    $para = Get-ContainingParagraph $rng
    $pRng = $para.Range
    $xml = $pkgHeader + $paragraphInnerXml + $pkgFooter
    $pRng.InsertXML($xml)
}

# 1) Command-line interface cell: "Type arg -f filter" -> split "arg" out
#    with spell-check proof-error markers around it.
Set-ParagraphXml ' arg -f filter' (
    '<w:p w14:paraId="66CBF6AC" w14:textId="445D5915" w:rsidR="004474EB" w:rsidRDefault="00345E99" w:rsidP="00345E99">' +
      '<w:pPr><w:pStyle w:val="HTMLPreformatted"/></w:pPr>' +
      '<w:r w:rsidRPr="00345E99">' +
        '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
        '<w:t>Type</w:t>' +
      '</w:r>' +
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>arg</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t xml:space="preserve"> -f filter</w:t></w:r>' +
    '</w:p>'
)

# 2) <diagnostics-general> table: turn each "Checking the ..." cell into a
#    screen-capture placeholder (blue #0070C0, underlined), matching the
#    style used elsewhere in the document for that convention.
$rPrCommon = '<w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:color w:val="0070C0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/>'

Set-ParagraphXml 'Checking the struts' (
    '<w:p w14:paraId="26140490" w14:textId="17C4A4A1" w:rsidR="002505C0" w:rsidRPr="00BD7996" w:rsidRDefault="00BD7996" w:rsidP="00913452">' +
      '<w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr>' + $rPrCommon + '</w:rPr></w:pPr>' +
      '<w:r w:rsidRPr="00BD7996"><w:rPr>' + $rPrCommon + '</w:rPr><w:t>Checking the struts</w:t></w:r>' +
    '</w:p>'
)

Set-ParagraphXml 'Checking the ball joints' (
    '<w:p w14:paraId="037FF1D0" w14:textId="740FEA4E" w:rsidR="002505C0" w:rsidRPr="00BD7996" w:rsidRDefault="00BD7996" w:rsidP="00913452">' +
      '<w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr>' + $rPrCommon + '</w:rPr></w:pPr>' +
      '<w:r><w:rPr>' + $rPrCommon + '</w:rPr><w:t>Checking the ball joints</w:t></w:r>' +
    '</w:p>'
)

Set-ParagraphXml 'Checking the exhaust' (
    '<w:p w14:paraId="5C4AB742" w14:textId="23E42DCB" w:rsidR="002505C0" w:rsidRPr="00BD7996" w:rsidRDefault="00BD7996" w:rsidP="00913452">' +
      '<w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr>' + $rPrCommon + '</w:rPr></w:pPr>' +
      '<w:r><w:rPr>' + $rPrCommon + '</w:rPr><w:t>Checking the exhaust</w:t></w:r>' +
    '</w:p>'
)

Set-ParagraphXml 'Checking the brakes' (
    '<w:p w14:paraId="19B8B259" w14:textId="57E9BBB0" w:rsidR="002505C0" w:rsidRPr="00BD7996" w:rsidRDefault="00BD7996" w:rsidP="00913452">' +
      '<w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr>' + $rPrCommon + '</w:rPr></w:pPr>' +
      '<w:r><w:rPr>' + $rPrCommon + '</w:rPr><w:t>Checking the brakes</w:t></w:r>' +
    '</w:p>'
)

Set-ParagraphXml 'Checking the CV joints' (
    '<w:p w14:paraId="7B91DE02" w14:textId="76B82B10" w:rsidR="002505C0" w:rsidRPr="00BD7996" w:rsidRDefault="00BD7996" w:rsidP="00913452">' +
      '<w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr>' + $rPrCommon + '</w:rPr></w:pPr>' +
      '<w:r><w:rPr>' + $rPrCommon + '</w:rPr><w:t>Checking the CV joints</w:t></w:r>' +
    '</w:p>'
)

Set-ParagraphXml 'Checking the catalytic converter' (
    '<w:p w14:paraId="772C3322" w14:textId="2C868C96" w:rsidR="002505C0" w:rsidRPr="00BD7996" w:rsidRDefault="00BD7996" w:rsidP="00913452">' +
      '<w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr>' + $rPrCommon + '</w:rPr></w:pPr>' +
      '<w:r><w:rPr>' + $rPrCommon + '</w:rPr><w:t>Checking the catalytic converter</w:t></w:r>' +
    '</w:p>'
)

Set-ParagraphXml 'Checking the wheels' (
    '<w:p w14:paraId="566B1C30" w14:textId="74ABEC2A" w:rsidR="002505C0" w:rsidRPr="00BD7996" w:rsidRDefault="00BD7996" w:rsidP="00913452">' +
      '<w:pPr><w:pStyle w:val="HTMLPreformatted"/><w:rPr>' + $rPrCommon + '</w:rPr></w:pPr>' +
      '<w:r><w:rPr>' + $rPrCommon + '</w:rPr><w:t>Checking the wheels</w:t></w:r>' +
    '</w:p>'
)

Write-Host "Done"
